# Apply "Add UK software/tech to list" edit to the Overview - Software workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Main" sheet (sheet2.xml)
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Cazoo (row 24): Sub-sector tag gains ", Cloud"
# ---------------------------------------------------------------------------
$ws.Cells.Item(24,35).Value = "E-Commerce, Logistics, Cloud"   # AI24

# ---------------------------------------------------------------------------
# 2. Auto Trader (row 36): stray backtick placed in column Q
# ---------------------------------------------------------------------------
$ws.Cells.Item(36,17).Value = [string]([char]96)   # Q36 = `

# ---------------------------------------------------------------------------
# 3. Insert 9 new rows right after row 40 (Trustpilot), pushing the FX-rate
#    table (previously rows 44-46) down to rows 53-55. Excel automatically
#    re-points the dependent formulas (F45 -> F54, $F$45 -> $F$54, G45 -> G54).
# ---------------------------------------------------------------------------
$ws.Rows("41:49").Insert(-4121)

# The insert duplicates row 40's per-cell formatting across the new rows;
# strip the cells that should stay empty/unformatted so the sheet matches
# the target layout.
$ws.Range("B41:F42").Clear()
$ws.Range("AD41:AH42").Clear()

$ws.Range("AD43:AF43").Clear()

$ws.Range("E44:E44").Clear()
$ws.Range("AD44:AF44").Clear()

$ws.Range("B45:F45").Clear()
$ws.Range("AD45:AH45").Clear()

$ws.Range("D46:F46").Clear()
$ws.Range("AD46:AH46").Clear()

$ws.Range("D47:F47").Clear()
$ws.Range("AD47:AH47").Clear()

$ws.Range("D48:F48").Clear()
$ws.Range("AD48:AH48").Clear()

$ws.Range("B49:F49").Clear()
$ws.Range("AD49:AH49").Clear()

# ---------------------------------------------------------------------------
# 4. New streaming rows (43 = Netflix, 44 = Spotify)
# ---------------------------------------------------------------------------
$ws.Cells.Item(44,2).Value = '$SPOT'
$ws.Cells.Item(43,2).Value = '$NFLX'
$ws.Cells.Item(44,3).Value = 'Spotify Technology S.A.'
$ws.Cells.Item(43,3).Value = 'Netflix, Inc.'
$ws.Cells.Item(43,35).Value = 'Movies'
$ws.Cells.Item(44,35).Value = 'Music'
$ws.Cells.Item(43,34).Value = 'SaaS Streaming'
$ws.Cells.Item(44,34).Value = 'SaaS Streaming'

$ws.Cells.Item(43,4).Value = 'NASDAQ'
$ws.Cells.Item(43,5).Value = '$'
$ws.Cells.Item(43,6).Formula = '=332.82*F54'
$ws.Cells.Item(43,6).NumberFormat = '#,##0.00'
$ws.Cells.Item(43,7).Value = 445.02
$ws.Cells.Item(43,8).Formula = '=G43*F43'

$ws.Cells.Item(44,4).Value = 'NYSE'
$ws.Cells.Item(44,6).Formula = '=92.06*F54'
$ws.Cells.Item(44,6).NumberFormat = '#,##0.00'
$ws.Cells.Item(44,7).Value = 193.13
$ws.Cells.Item(44,8).Formula = '=G44*F44'

# ---------------------------------------------------------------------------
# 5. New UK software/tech rows (46 = Idox, 47 = Computacenter, 48 = Softcat)
# ---------------------------------------------------------------------------
$ws.Cells.Item(47,3).Value = 'Computacenter Plc'
$ws.Cells.Item(46,3).Value = 'Idox Plc'
$ws.Cells.Item(48,3).Value = 'Softcat Plc'
$ws.Cells.Item(46,2).Value = '£IDOX'
$ws.Cells.Item(47,2).Value = '£CCC'
$ws.Cells.Item(48,2).Value = '£SCT'

# ---------------------------------------------------------------------------
# 6. FX-rate table (now rows 53-55): USDGBP rate update 0.83 -> 0.82
# ---------------------------------------------------------------------------
$ws.Cells.Item(54,6).Value = 0.82

# ---------------------------------------------------------------------------
# 7. View state: selected cell / scroll position
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("D48").Select()

Write-Host "edit applied"
